$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 entirely (its data is being dropped; row 4 becomes the last data row)
$ws.Rows.Item(5).Delete()

# Update row 4 with the merged/updated values
$ws.Range("A4").Value = "25/08/2016"
$ws.Range("B4").Value = "Done"
$ws.Range("C4").Value = "Done"
$ws.Range("D4").Value = "Done"
$ws.Range("E4").Value = "Done"
$ws.Range("F4").Value = "Done"
$ws.Range("G4").Value = "Done"
$ws.Range("H4").Value = "In progress"
$ws.Range("I4").Value = "In progress"
$ws.Range("J4").Value = "TD"
$ws.Range("K4").Value = "TD"
$ws.Range("L4").Value = "TD"

# Update sheet view selection to match the post-edit state
$ws.Range("B5").Select()
